$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 759
$ws.Range("F6").Value = 2432
$ws.Range("F8").Value = 1794
$ws.Range("F9").Value = 3055
$ws.Range("F10").Value = 183
$ws.Range("F11").Value = 4537
$ws.Range("F12").Value = 407
$ws.Range("F13").Value = 231
$ws.Range("F14").Value = 142
$ws.Range("F15").Value = 573
$ws.Range("F16").Value = 269
$ws.Range("F17").Value = 621
$ws.Range("F20").Value = 119
$ws.Range("F22").Value = 4567
$ws.Range("F24").Value = 4106
$ws.Range("F25").Value = 1150
$ws.Range("F26").Value = 221
$ws.Range("F27").Value = 599
$ws.Range("F29").Value = 97
$ws.Range("F30").Value = 667
$ws.Range("F31").Value = 627
$ws.Range("F32").Value = 609
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 36
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 212
$ws.Range("F4").Value = 24
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 212
$ws.Range("F5").Value = 24
$ws.Range("F8").Value = 759
$ws.Range("F9").Value = 2432
$ws.Range("F11").Value = 1794
$ws.Range("F13").Value = 3055
$ws.Range("F14").Value = 183
$ws.Range("F15").Value = 4537
$ws.Range("F16").Value = 407
$ws.Range("F17").Value = 231
$ws.Range("F18").Value = 142
$ws.Range("F19").Value = 573
$ws.Range("F20").Value = 269
$ws.Range("F21").Value = 621
$ws.Range("F25").Value = 119
$ws.Range("F27").Value = 4567
$ws.Range("F29").Value = 4106
$ws.Range("F30").Value = 1150
$ws.Range("F31").Value = 221
$ws.Range("F32").Value = 599
$ws.Range("F35").Value = 667
$ws.Range("F36").Value = 627
$ws.Range("F37").Value = 609
$ws.Range("F39").Value = 36
